$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.753.55'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +8.19%  '
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.405.73'
$ws.Range('D3').Style = $origStyle
$ws.Range('E3').Value = '  +5.16%  '
$ws.Range('E4').Value = '  +0.05%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '419.21'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +6.32%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.65'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +8.16%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.599'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +6.30%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +5.39%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.115'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +18.11%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.80'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +4.99%  '
$ws.Range('E12').Value = '  +1.40%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.924.27'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +4.82%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.58'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +5.58%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.12'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +6.36%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.428.79'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +5.85%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.06'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +2.45%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.564.88'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +8.31%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.92'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000119'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +13.16%  '
$ws.Range('B21').Value = 'ImmutableX'
$ws.Range('C21').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.43'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +2.91%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.24'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +2.27%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '306.81'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  +2.93%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '76.42'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +3.87%  '
$ws.Range('E25').Value = '  +3.60%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.47'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +5.90%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.50'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +2.60%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.97'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +3.41%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.69'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +6.34%  '
$ws.Range('E30').Value = '  +6.84%  '
$ws.Range('E31').Value = '  +6.88%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.59'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +6.05%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.56'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +21.64%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +0.08%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.03'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +7.52%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0513'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +6.20%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.63'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +1.92%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.14'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +1.34%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -0.07%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.46'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -1.43%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.67'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +2.72%  '
$ws.Range('E42').Value = '  +4.41%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.95'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.05'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.292'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +3.91%  '
$ws.Range('E46').Value = '  +1.38%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.32'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +12.39%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.17'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +5.90%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.178.61'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('E51').Value = '  -0.79%  '
